$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 27; this shifts existing rows 27..106 down to 28..107
$ws.Rows.Item(27).Insert()

# The inherited/static columns are the same for every data row in this sheet
$ws.Cells.Item(27, 1).Value = 7
$ws.Cells.Item(27, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(27, 3).Value = "Ñuble"
$ws.Cells.Item(27, 4).Value = 44914
$ws.Cells.Item(27, 5).Value = 16
$ws.Cells.Item(27, 6).Value = 100112031
$ws.Cells.Item(27, 7).Value = "Poroto verde"
$ws.Cells.Item(27, 8).Value = "Sin especificar"
$ws.Cells.Item(27, 9).Value = "Primera"
$ws.Cells.Item(27, 10).Value = 50
$ws.Cells.Item(27, 11).Value = 30000
$ws.Cells.Item(27, 12).Value = 30000
$ws.Cells.Item(27, 13).Value = 30000
$ws.Cells.Item(27, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(27, 15).Value = "Región del Maule"
$ws.Cells.Item(27, 16).Value = 1200
$ws.Cells.Item(27, 17).Value = 25
$ws.Cells.Item(27, 18).Value = "Hortaliza"

# Match the date-cell style used by the rest of column D
$ws.Cells.Item(27, 4).NumberFormat = $ws.Cells.Item(28, 4).NumberFormat
